# "Mails für die Reisen zugefügt" - update Toskana 2017 calculation sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rabatt Nichtradf. (B13): 7 -> 6. All dependent formulas (I18, I19, I21,
# I22, B25, B26, I26, B27, ...) recalc automatically from this single
# input change.
$ws.Range("B13").Value = 6

# Move the active selection to B14, matching the sheet's saved cursor
# position.
$ws.Range("B14").Select()
